$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column A: EmpUnqID value change
$ws.Range("A2").Value = 102434

# Row 2: every "GN" data value becomes "G1" (a rename - "GN" shared string
# is fully replaced by a new "G1" shared string), except column D which is
# switched to "WO" (mirroring the "WO" that used to sit in column C).
$ws.Range("B2").Value = "G1"
$ws.Range("C2").Value = "G1"
$ws.Range("D2").Value = "WO"
$ws.Range("E2").Value = "G1"
$ws.Range("F2").Value = "G1"
$ws.Range("G2").Value = "G1"
$ws.Range("H2").Value = "G1"
$ws.Range("I2").Value = "G1"
$ws.Range("J2").Value = "G1"
$ws.Range("K2").Value = "WO"
$ws.Range("L2").Value = "G1"
$ws.Range("M2").Value = "G1"
$ws.Range("N2").Value = "G1"
$ws.Range("O2").Value = "G1"
$ws.Range("P2").Value = "G1"
$ws.Range("Q2").Value = "G1"
$ws.Range("R2").Value = "WO"
$ws.Range("S2").Value = "G1"
$ws.Range("T2").Value = "G1"
$ws.Range("U2").Value = "G1"
$ws.Range("V2").Value = "G1"
$ws.Range("W2").Value = "G1"
$ws.Range("X2").Value = "G1"
$ws.Range("Y2").Value = "WO"
$ws.Range("Z2").Value = "G1"
$ws.Range("AA2").Value = "G1"
$ws.Range("AB2").Value = "G1"
$ws.Range("AC2").Value = "G1"
$ws.Range("AD2").Value = "G1"
$ws.Range("AE2").Value = "G1"
$ws.Range("AF2").Value = "WO"

# AC2 picks up the same (centered) formatting that its neighbours already
# use, instead of the "no alignment" style it had before.
$ws.Range("AB2").Copy()
$ws.Range("AC2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active cell of the bottom-right pane from J4 to E5.
[void]$ws.Range("E5").Select()
